$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1792.625
$ws.Range("J112").Value = 1792.625
$ws.Range("L112").Value = 5377.875
$ws.Range("N112").Value = -7593.875
$ws.Range("H129").Value = 1478.3556
$ws.Range("I129").Value = 1939.7142
$ws.Range("J129").Value = 1393.3684
$ws.Range("K129").Value = 5819.142599999999
$ws.Range("L129").Value = 4180.1052
$ws.Range("M129").Value = -819.1425999999992
$ws.Range("N129").Value = -14180.1052
$ws.Range("H137").Value = 3436.4517
$ws.Range("I137").Value = 1020
$ws.Range("J137").Value = 3901.1538
$ws.Range("K137").Value = 3060
$ws.Range("L137").Value = 11703.4614
$ws.Range("M137").Value = -510
$ws.Range("N137").Value = -16803.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4870.033
$ws.Range("I2").Value = 5572.08
$ws.Range("J2").Value = 1359.8
$ws.Range("K2").Value = 5572.08
$ws.Range("L2").Value = 1359.8
$ws.Range("M2").Value = -5459.08
$ws.Range("N2").Value = -1585.8
$ws.Range("H61").Value = 4578.65
$ws.Range("I61").Value = 2999.75
$ws.Range("J61").Value = 4973.375
$ws.Range("K61").Value = 2999.75
$ws.Range("L61").Value = 4973.375
$ws.Range("M61").Value = -2787.75
$ws.Range("N61").Value = -5397.375
$ws.Range("H97").Value = 954
$ws.Range("I97").Value = 731.4286
$ws.Range("J97").Value = 1148.75
$ws.Range("K97").Value = 731.4286
$ws.Range("L97").Value = 1148.75
$ws.Range("M97").Value = -235.4286
$ws.Range("N97").Value = -2140.75
$ws.Range("H116").Value = 4870.033
$ws.Range("I116").Value = 5572.08
$ws.Range("J116").Value = 1359.8
$ws.Range("K116").Value = 5572.08
$ws.Range("L116").Value = 1359.8
$ws.Range("M116").Value = -3278.08
$ws.Range("N116").Value = -5947.8
$ws.Range("H136").Value = 4578.65
$ws.Range("I136").Value = 2999.75
$ws.Range("J136").Value = 4973.375
$ws.Range("K136").Value = 8999.25
$ws.Range("L136").Value = 14920.125
$ws.Range("M136").Value = -6449.25
$ws.Range("N136").Value = -20020.125
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4870.033
$ws.Range("I3").Value = 5572.08
$ws.Range("J3").Value = 1359.8
$ws.Range("K3").Value = 5572.08
$ws.Range("L3").Value = 1359.8
$ws.Range("M3").Value = -5458.08
$ws.Range("N3").Value = -1587.8
$ws.Range("H105").Value = 1908.7916
$ws.Range("I105").Value = 1927.1666
$ws.Range("J105").Value = 1853.6666
$ws.Range("K105").Value = 1927.1666
$ws.Range("L105").Value = 1853.6666
$ws.Range("M105").Value = -180.1666
$ws.Range("N105").Value = -5347.6666
$ws.Range("H124").Value = 50996
$ws.Range("J124").Value = 50996
$ws.Range("L124").Value = 50996
$ws.Range("N124").Value = -60816
$ws.Range("H134").Value = 4310.519
$ws.Range("I134").Value = 2944.8667
$ws.Range("J134").Value = 4864.162
$ws.Range("K134").Value = 8834.6001
$ws.Range("L134").Value = 14592.486
$ws.Range("M134").Value = -6299.6001
$ws.Range("N134").Value = -19662.486

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H111").Value = 46996
$ws.Range("J111").Value = 46996
$ws.Range("L111").Value = 46996
$ws.Range("N111").Value = -55176
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 503.25
$ws.Range("I13").Value = 455.5
$ws.Range("J13").Value = 551
$ws.Range("K13").Value = 1366.5
$ws.Range("L13").Value = 1653
$ws.Range("M13").Value = -1198.5
$ws.Range("N13").Value = -1989
$ws.Range("H68").Value = 576.7692
$ws.Range("I68").Value = 442.57144
$ws.Range("J68").Value = 733.3333
$ws.Range("K68").Value = 1327.71432
$ws.Range("L68").Value = 2199.9999
$ws.Range("M68").Value = -516.71432
$ws.Range("N68").Value = -3821.9999
$ws.Range("H71").Value = 576.7692
$ws.Range("I71").Value = 442.57144
$ws.Range("J71").Value = 733.3333
$ws.Range("K71").Value = 3983.14296
$ws.Range("L71").Value = 6599.9997
$ws.Range("M71").Value = 72.85703999999987
$ws.Range("N71").Value = -14711.9997
$ws.Range("H82").Value = 3496.5
$ws.Range("I82").Value = 3496.5
$ws.Range("K82").Value = 10489.5
$ws.Range("M82").Value = -10083.5
$ws.Range("H85").Value = 3496.5
$ws.Range("I85").Value = 3496.5
$ws.Range("K85").Value = 10489.5
$ws.Range("M85").Value = -9085.5
$ws.Range("H86").Value = 1572.375
$ws.Range("I86").Value = 790
$ws.Range("J86").Value = 1833.1666
$ws.Range("K86").Value = 2370
$ws.Range("L86").Value = 5499.4998
$ws.Range("M86").Value = -1184
$ws.Range("N86").Value = -7871.4998
$ws.Range("H89").Value = 1572.375
$ws.Range("I89").Value = 790
$ws.Range("J89").Value = 1833.1666
$ws.Range("K89").Value = 7110
$ws.Range("L89").Value = 16498.4994
$ws.Range("M89").Value = -1182
$ws.Range("N89").Value = -28354.4994
$ws.Range("H132").Value = 1435.45
$ws.Range("I132").Value = 1060.4
$ws.Range("J132").Value = 1810.5
$ws.Range("K132").Value = 9543.6
$ws.Range("L132").Value = 16294.5
$ws.Range("M132").Value = -7013.6
$ws.Range("N132").Value = -21354.5
$ws.Range("H134").Value = 33370718
$ws.Range("I134").Value = 55613676
$ws.Range("J134").Value = 6280.6665
$ws.Range("K134").Value = 166841028
$ws.Range("L134").Value = 18841.9995
$ws.Range("M134").Value = -166835958
$ws.Range("N134").Value = -28981.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3829.587
$ws.Range("I97").Value = 1440.6154
$ws.Range("J97").Value = 6935.25
$ws.Range("K97").Value = 1440.6154
$ws.Range("L97").Value = 6935.25
$ws.Range("M97").Value = -944.6153999999999
$ws.Range("N97").Value = -7927.25
$ws.Range("H122").Value = 2218.9412
$ws.Range("I122").Value = 1883.091
$ws.Range("J122").Value = 2834.6667
$ws.Range("K122").Value = 5649.272999999999
$ws.Range("L122").Value = 8504.000100000001
$ws.Range("M122").Value = -3199.272999999999
$ws.Range("N122").Value = -13404.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 23955
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 23955
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 23955
$ws.Range("N33").Value = -24455
$ws.Range("H36").Value = 23955
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 23955
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 23955
$ws.Range("N36").Value = -24455
$ws.Range("H92").Value = 21433.334
$ws.Range("J92").Value = 21433.334
$ws.Range("L92").Value = 21433.334
$ws.Range("N92").Value = -26425.334
$ws.Range("H139").Value = 45337
$ws.Range("J139").Value = 45337
$ws.Range("L139").Value = 45337
$ws.Range("N139").Value = -55617
$ws.Range("M33").ClearContents()
$ws.Range("M36").ClearContents()
